$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ROLL NO"
$ws.Range("D1").Value = "FATHER NAME"
$ws.Range("F1").Value = "COURSE ID"

$ws.Range("D7").Select()
